# "completed part 3 and 4" — this worksheet held a leftover "Audio"
# column (the header cell + the "audio.mp3" value merged across the
# three question rows) that is no longer needed now that part 3/4 are
# finished. Remove the whole column, which shifts every other column
# left by one, drops the now-empty merge, and shrinks the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the column first (mirrors clicking the "A" column header)
# before deleting it, same as a user would do interactively.
[void]$ws.Columns("A").Select()
$ws.Columns("A").Delete() | Out-Null
